$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "BK"
$ws.Range("B13").Value = "PM"
$ws.Range("C13").Value = 1234
$ws.Range("D13").Value = "finished"
$ws.Range("E13").Value = "cdf24541-9127-458c-9cd5-b8a12f4b929c"
$ws.Range("F13").Value = "2025-03-11T13:53:26.610680"
$ws.Range("G13").Value = "2025-03-11T13:55:06.282262"

$ws.Range("A14").Value = "BK"
$ws.Range("B14").Value = "PM"
$ws.Range("C14").Value = 2345
$ws.Range("D14").Value = "finished"
$ws.Range("E14").Value = "6d8fbad6-85c0-425a-a3dc-61f9817d2009"
$ws.Range("F14").Value = "2025-03-11T13:53:30.710633"
$ws.Range("G14").Value = "2025-03-11T13:53:32.391061"

$ws.Range("A15").Value = "BK"
$ws.Range("B15").Value = "PM"
$ws.Range("C15").Value = 1234
$ws.Range("D15").Value = "finished"
$ws.Range("E15").Value = "6bfd6e1c-2f99-4993-94fb-c8f296212a53"
$ws.Range("F15").Value = "2025-03-11T13:54:54.693445"
$ws.Range("G15").Value = "2025-03-11T14:08:24.249050"

$ws.Range("A16").Value = "BK"
$ws.Range("B16").Value = "PM"
$ws.Range("C16").Value = 2456
$ws.Range("D16").Value = "finished"
$ws.Range("E16").Value = "27daef15-51da-4a46-a80d-325cc956802e"
$ws.Range("F16").Value = "2025-03-11T13:54:59.197709"
$ws.Range("G16").Value = "2025-03-11T14:08:24.386535"

$ws.Range("A17").Value = "BK"
$ws.Range("B17").Value = "PM"
$ws.Range("C17").Value = 1234
$ws.Range("D17").Value = "finished"
$ws.Range("E17").Value = "31388b15-ac05-400d-a303-d5f639f45404"
$ws.Range("F17").Value = "2025-03-11T13:55:40.379223"
$ws.Range("G17").Value = "2025-03-12T11:03:48.781585"

$ws.Range("A18").Value = "BK"
$ws.Range("B18").Value = "PM"
$ws.Range("C18").Value = 1245
$ws.Range("D18").Value = "finished"
$ws.Range("E18").Value = "3864c7a9-753a-467e-9e3d-d05048ff16d7"
$ws.Range("F18").Value = "2025-03-11T13:55:44.585014"
$ws.Range("G18").Value = "2025-03-12T11:03:45.695443"

$ws.Range("A19").Value = "BK"
$ws.Range("B19").Value = "PM"
$ws.Range("C19").Value = 5678
$ws.Range("D19").Value = "finished"
$ws.Range("E19").Value = "e2b508a8-f1a7-42e2-a070-a2f670779e15"
$ws.Range("F19").Value = "2025-03-11T13:55:47.751952"
$ws.Range("G19").Value = "2025-03-11T14:09:28.875654"

$ws.Range("A20").Value = "BK"
$ws.Range("B20").Value = "PM"
$ws.Range("C20").Value = 1234
$ws.Range("D20").Value = "finished"
$ws.Range("E20").Value = "1693ac31-331b-4156-90d9-98c1b862bfdc"
$ws.Range("F20").Value = "2025-03-11T13:59:02.371825"
$ws.Range("G20").Value = "2025-03-12T11:03:47.019402"

$ws.Range("A21").Value = "BK"
$ws.Range("B21").Value = "PM"
$ws.Range("C21").Value = 1234
$ws.Range("D21").Value = "finished"
$ws.Range("E21").Value = "c6559d3d-b592-4d32-b93a-d6f4ebaf2da4"
$ws.Range("F21").Value = "2025-03-12T08:09:09.694421"
$ws.Range("G21").Value = "2025-03-12T11:03:47.790109"

$ws.Range("A22").Value = "BK"
$ws.Range("B22").Value = "PM"
$ws.Range("C22").Value = 1233
$ws.Range("D22").Value = "finished"
$ws.Range("E22").Value = "cbf07666-709f-480c-ac5c-0b2b784cc3b9"
$ws.Range("F22").Value = "2025-03-12T08:09:13.875922"
$ws.Range("G22").Value = "2025-03-12T11:03:49.283978"

$ws.Range("A23").Value = "BK"
$ws.Range("B23").Value = "PM"
$ws.Range("C23").Value = 3456
$ws.Range("D23").Value = "finished"
$ws.Range("E23").Value = "fa885632-045e-48ce-bca3-39f8a6d56b28"
$ws.Range("F23").Value = "2025-03-12T08:25:51.334453"
$ws.Range("G23").Value = "2025-03-12T09:42:06.423077"

$ws.Range("A24").Value = "BK"
$ws.Range("B24").Value = "PM"
$ws.Range("C24").Value = 1234
$ws.Range("D24").Value = "finished"
$ws.Range("E24").Value = "9106c101-5720-4d4a-a8d0-337c6d325549"
$ws.Range("F24").Value = "2025-03-12T08:27:38.952059"
$ws.Range("G24").Value = "2025-03-12T08:27:40.743786"
Write-Output "Rows added."
